$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.04741066666666666
$ws.Range("H2").Value = 0.142232
$ws.Range("I2").Value = 0.003188134523263584
$ws.Range("J2").Value = 0.003188134523263585
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.040495
$ws.Range("N2").Value = 0.121485
$ws.Range("O2").Value = 0.002191743187342868
$ws.Range("P2").Value = 0.002191743187342869
$ws.Range("Q2").Value = 0.001919894946666666
$ws.Range("R2").Value = 0.01727905452
$ws.Range("S2").Value = 0.000006987572121695564
$ws.Range("T2").Value = 0.000006987572121695566

# Row 3
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.04741066666666666
$ws.Range("H3").Value = 0.142232
$ws.Range("I3").Value = 0.003188134523263584
$ws.Range("J3").Value = 0.003188134523263585
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 16.98312366666667
$ws.Range("N3").Value = 50.949371
$ws.Range("O3").Value = 0.9191911494312409
$ws.Range("P3").Value = 0.9191911494312409
$ws.Range("Q3").Value = 0.8051812151191111
$ws.Range("R3").Value = 7.246630936072
$ws.Range("S3").Value = 0.002930505036980076
$ws.Range("T3").Value = 0.002930505036980076

# Row 4
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.04741066666666666
$ws.Range("H4").Value = 0.142232
$ws.Range("I4").Value = 0.003188134523263584
$ws.Range("J4").Value = 0.003188134523263585
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 1.452542333333333
$ws.Range("N4").Value = 4.357627
$ws.Range("O4").Value = 0.07861710738141615
$ws.Range("P4").Value = 0.07861710738141615
$ws.Range("Q4").Value = 0.06886600038488888
$ws.Range("R4").Value = 0.6197940034639999
$ws.Range("S4").Value = 0.0002506419141618132
$ws.Range("T4").Value = 0.0002506419141618133

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 0.237305
$ws.Range("H5").Value = 0.711915
$ws.Range("I5").Value = 0.01595759596384214
$ws.Range("J5").Value = 0.01595759596384214
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.040495
$ws.Range("N5").Value = 0.121485
$ws.Range("O5").Value = 0.002191743187342868
$ws.Range("P5").Value = 0.002191743187342869
$ws.Range("Q5").Value = 0.009609665974999999
$ws.Range("R5").Value = 0.08648699377499999
$ws.Range("S5").Value = 0.00003497495224012105
$ws.Range("T5").Value = 0.00003497495224012106

# Row 6
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 0.237305
$ws.Range("H6").Value = 0.711915
$ws.Range("I6").Value = 0.01595759596384214
$ws.Range("J6").Value = 0.01595759596384214
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 16.98312366666667
$ws.Range("N6").Value = 50.949371
$ws.Range("O6").Value = 0.9191911494312409
$ws.Range("P6").Value = 0.9191911494312409
$ws.Range("Q6").Value = 4.030180161718333
$ws.Range("R6").Value = 36.271621455465
$ws.Range("S6").Value = 0.01466808097616338
$ws.Range("T6").Value = 0.01466808097616338

# Row 7
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 0.237305
$ws.Range("H7").Value = 0.711915
$ws.Range("I7").Value = 0.01595759596384214
$ws.Range("J7").Value = 0.01595759596384214
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 1.452542333333333
$ws.Range("N7").Value = 4.357627
$ws.Range("O7").Value = 0.07861710738141615
$ws.Range("P7").Value = 0.07861710738141615
$ws.Range("Q7").Value = 0.3446955584116667
$ws.Range("R7").Value = 3.102260025705
$ws.Range("S7").Value = 0.00125454003543863
$ws.Range("T7").Value = 0.00125454003543863

# Row 8
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 14.58625866666667
$ws.Range("H8").Value = 43.758776
$ws.Range("I8").Value = 0.9808542695128942
$ws.Range("J8").Value = 0.9808542695128943
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 0.040495
$ws.Range("N8").Value = 0.121485
$ws.Range("O8").Value = 0.002191743187342868
$ws.Range("P8").Value = 0.002191743187342869
$ws.Range("Q8").Value = 0.5906705447066666
$ws.Range("R8").Value = 5.316034902359999
$ws.Range("S8").Value = 0.002149780662981052
$ws.Range("T8").Value = 0.002149780662981052

# Row 9
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 14.58625866666667
$ws.Range("H9").Value = 43.758776
$ws.Range("I9").Value = 0.9808542695128942
$ws.Range("J9").Value = 0.9808542695128943
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 16.98312366666667
$ws.Range("N9").Value = 50.949371
$ws.Range("O9").Value = 0.9191911494312409
$ws.Range("P9").Value = 0.9191911494312409
$ws.Range("Q9").Value = 247.7202347699885
$ws.Range("R9").Value = 2229.482112929896
$ws.Range("S9").Value = 0.9015925634180975
$ws.Range("T9").Value = 0.9015925634180975

# Row 10
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 14.58625866666667
$ws.Range("H10").Value = 43.758776
$ws.Range("I10").Value = 0.9808542695128942
$ws.Range("J10").Value = 0.9808542695128943
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 1.452542333333333
$ws.Range("N10").Value = 4.357627
$ws.Range("O10").Value = 0.07861710738141615
$ws.Range("P10").Value = 0.07861710738141615
$ws.Range("Q10").Value = 21.18715819828356
$ws.Range("R10").Value = 190.684423784552
$ws.Range("S10").Value = 0.0771119254318157
$ws.Range("T10").Value = 0.07711192543181572
